$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.590.49'
$ws.Range('D3').Value = '1.698.33'
$ws.Range('E3').Value = '  +2.13%  '
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'315.21"
$ws.Range('E5').Value = '  +1.81%  '
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = "'0.3944"
$ws.Range('E7').Value = '  +1.16%  '
$ws.Range('E8').Value = '  +1.62%  '
$ws.Range('D9').Value = "'1.549"
$ws.Range('E9').Value = '  +8.14%  '
$ws.Range('D10').Value = "'55.29"
$ws.Range('E10').Value = '  +14.41%  '
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('D12').Value = "'0.08801"
$ws.Range('E12').Value = '  +1.90%  '
$ws.Range('D13').Value = "'7.306"
$ws.Range('E13').Value = '  +11.80%  '
$ws.Range('D14').Value = "'23.36"
$ws.Range('E14').Value = '  +2.79%  '
$ws.Range('D15').Value = "'0.00001332"
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').Value = "'7.635"
$ws.Range('E16').Value = '  +5.97%  '
$ws.Range('D17').Value = '1.697.67'
$ws.Range('E17').Value = '  +2.32%  '
$ws.Range('D18').Value = "'101.19"
$ws.Range('E18').Value = '  +0.79%  '
$ws.Range('D19').Value = "'0.07030"
$ws.Range('E19').Value = '  +3.76%  '
$ws.Range('D20').Value = "'19.84"
$ws.Range('E20').Value = '  +4.34%  '
$ws.Range('D21').Value = "'6.938"
$ws.Range('E21').Value = '  +4.54%  '
$ws.Range('D22').Value = "'1.002"
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('D24').Value = '24.578.72'
$ws.Range('E24').Value = '  +3.23%  '
$ws.Range('D25').Value = "'2.978"
$ws.Range('E25').Value = '  +8.45%  '
$ws.Range('D26').Value = "'2.340"
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('D27').Value = "'22.41"
$ws.Range('E27').Value = '  +2.98%  '
$ws.Range('D28').Value = "'160.06"
$ws.Range('E28').Value = '  +1.75%  '
$ws.Range('D29').Value = "'5.246"
$ws.Range('E29').Value = '  +1.95%  '
$ws.Range('D30').Value = "'134.01"
$ws.Range('E30').Value = '  +3.28%  '
$ws.Range('D31').Value = "'7.637"
$ws.Range('E31').Value = '  +21.06%  '
$ws.Range('D32').Value = "'1.115"
$ws.Range('E32').Value = '  -2.49%  '
$ws.Range('D33').Value = '1.883.11'
$ws.Range('E33').Value = '  +2.25%  '
$ws.Range('D34').Value = "'7.478"
$ws.Range('E34').Value = '  +13.90%  '
$ws.Range('D35').Value = "'0.08579"
$ws.Range('E35').Value = '  -0.82%  '
$ws.Range('D36').Value = "'11.22"
$ws.Range('E36').Value = '  +8.67%  '
$ws.Range('D37').Value = "'1.980"
$ws.Range('E37').Value = '  +2.54%  '
$ws.Range('D38').Value = "'0.2759"
$ws.Range('E38').Value = '  +4.57%  '
$ws.Range('D39').Value = "'14.82"
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('D40').Value = "'0.02783"
$ws.Range('E40').Value = '  +10.59%  '
$ws.Range('D41').Value = "'0.09042"
$ws.Range('E41').Value = '  +2.95%  '
$ws.Range('D42').Value = "'1.474"
$ws.Range('E42').Value = '  +2.84%  '
$ws.Range('D43').Value = "'0.7767"
$ws.Range('E43').Value = '  +2.42%  '
$ws.Range('D44').Value = "'0.7298"
$ws.Range('E44').Value = '  +3.56%  '
$ws.Range('D45').Value = "'15.58"
$ws.Range('E45').Value = '  +4.71%  '
$ws.Range('D46').Value = "'2.514"
$ws.Range('E46').Value = '  +5.79%  '
$ws.Range('D47').Value = "'4.197"
$ws.Range('E47').Value = '  +2.83%  '
$ws.Range('D48').Value = "'1.000"
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('D49').Value = "'141.58"
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('D50').Value = "'1.299"
$ws.Range('E50').Value = '  +14.13%  '
$ws.Range('D51').Value = "'0.08025"
$ws.Range('E51').Value = '  +3.38%  '
